# Apply updated 'F' column (view/heat count) values per sheet, per commit diff
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 277  # was 276
$ws.Range("F4").Value = 970  # was 969
$ws.Range("F5").Value = 246  # was 245
$ws.Range("F7").Value = 667  # was 666
$ws.Range("F11").Value = 381  # was 379
$ws.Range("F12").Value = 180  # was 179
$ws.Range("F13").Value = 41  # was 38
$ws.Range("F14").Value = 768  # was 767
$ws.Range("F16").Value = 1908  # was 1907
$ws.Range("F17").Value = 423  # was 422
$ws.Range("F18").Value = 5858  # was 5806
$ws.Range("F20").Value = 510  # was 509
$ws.Range("F24").Value = 185  # was 184

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 499  # was 497
$ws.Range("F14").Value = 49  # was 48

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5429  # was 5424

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 5429  # was 5424
$ws.Range("F7").Value = 277  # was 276
$ws.Range("F12").Value = 499  # was 497
$ws.Range("F13").Value = 499  # was 497
$ws.Range("F14").Value = 970  # was 969
$ws.Range("F17").Value = 246  # was 245
$ws.Range("F19").Value = 667  # was 666
$ws.Range("F24").Value = 381  # was 379
$ws.Range("F25").Value = 180  # was 179
$ws.Range("F27").Value = 41  # was 38
$ws.Range("F29").Value = 768  # was 767
$ws.Range("F32").Value = 1908  # was 1907
$ws.Range("F33").Value = 423  # was 422
$ws.Range("F34").Value = 5858  # was 5806
$ws.Range("F35").Value = 49  # was 48
$ws.Range("F37").Value = 510  # was 509
$ws.Range("F42").Value = 185  # was 184
